$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 16.962335
$ws.Range("H2").Value = 50.887005
$ws.Range("I2").Value = 0.725422686224818
$ws.Range("J2").Value = 0.725422686224818
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.07605
$ws.Range("N2").Value = 18.22815
$ws.Range("O2").Value = 0.0302610603580868
$ws.Range("P2").Value = 0.0302610603580868
$ws.Range("Q2").Value = 103.06399557675
$ws.Range("R2").Value = 927.57596019075
$ws.Range("S2").Value = 0.02195205969297468
$ws.Range("T2").Value = 0.02195205969297468

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 16.962335
$ws.Range("H3").Value = 50.887005
$ws.Range("I3").Value = 0.725422686224818
$ws.Range("J3").Value = 0.725422686224818
$ws.Range("O3").Value = 0.2994824511432495
$ws.Range("P3").Value = 0.2994824511432494
$ws.Range("Q3").Value = 1019.986003619787
$ws.Range("R3").Value = 9179.874032578082
$ws.Range("S3").Value = 0.2172513641855288
$ws.Range("T3").Value = 0.2172513641855288

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 16.962335
$ws.Range("H4").Value = 50.887005
$ws.Range("I4").Value = 0.725422686224818
$ws.Range("J4").Value = 0.725422686224818
$ws.Range("M4").Value = 134.5792873333333
$ws.Range("N4").Value = 403.737862
$ws.Range("O4").Value = 0.6702564884986638
$ws.Range("P4").Value = 0.6702564884986637
$ws.Range("Q4").Value = 2282.778955809256
$ws.Range("R4").Value = 20545.01060228331
$ws.Range("S4").Value = 0.4862192623463145
$ws.Range("T4").Value = 0.4862192623463145

# Row 5
$ws.Range("H5").Value = 9.977416
$ws.Range("I5").Value = 0.1422336393407802
$ws.Range("J5").Value = 0.1422336393407802
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.07605
$ws.Range("N5").Value = 18.22815
$ws.Range("O5").Value = 0.0302610603580868
$ws.Range("P5").Value = 0.0302610603580868
$ws.Range("Q5").Value = 20.2077594956
$ws.Range("R5").Value = 181.8698354604
$ws.Range("S5").Value = 0.004304140745041699
$ws.Range("T5").Value = 0.004304140745041697

# Row 6
$ws.Range("H6").Value = 9.977416
$ws.Range("I6").Value = 0.1422336393407802
$ws.Range("J6").Value = 0.1422336393407802
$ws.Range("O6").Value = 0.2994824511432495
$ws.Range("P6").Value = 0.2994824511432494
$ws.Range("S6").Value = 0.04259647894480177
$ws.Range("T6").Value = 0.04259647894480176

# Row 7
$ws.Range("H7").Value = 9.977416
$ws.Range("I7").Value = 0.1422336393407802
$ws.Range("J7").Value = 0.1422336393407802
$ws.Range("M7").Value = 134.5792873333333
$ws.Range("N7").Value = 403.737862
$ws.Range("O7").Value = 0.6702564884986638
$ws.Range("P7").Value = 0.6702564884986637
$ws.Range("Q7").Value = 447.5845115693991
$ws.Range("R7").Value = 4028.260604124592
$ws.Range("S7").Value = 0.09533301965093674
$ws.Range("T7").Value = 0.09533301965093671

# Row 8
$ws.Range("G8").Value = 3.094551333333333
$ws.Range("H8").Value = 9.283654
$ws.Range("I8").Value = 0.1323436744344018
$ws.Range("J8").Value = 0.1323436744344018
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.07605
$ws.Range("N8").Value = 18.22815
$ws.Range("O8").Value = 0.0302610603580868
$ws.Range("P8").Value = 0.0302610603580868
$ws.Range("Q8").Value = 18.8026486289
$ws.Range("R8").Value = 169.2238376601
$ws.Range("S8").Value = 0.004004859920070422
$ws.Range("T8").Value = 0.004004859920070421

# Row 9
$ws.Range("G9").Value = 3.094551333333333
$ws.Range("H9").Value = 9.283654
$ws.Range("I9").Value = 0.1323436744344018
$ws.Range("J9").Value = 0.1323436744344018
$ws.Range("O9").Value = 0.2994824511432495
$ws.Range("P9").Value = 0.2994824511432494
$ws.Range("Q9").Value = 186.0828151008071
$ws.Range("R9").Value = 1674.745335907264
$ws.Range("S9").Value = 0.03963460801291886
$ws.Range("T9").Value = 0.03963460801291883

# Row 10
$ws.Range("G10").Value = 3.094551333333333
$ws.Range("H10").Value = 9.283654
$ws.Range("I10").Value = 0.1323436744344018
$ws.Range("J10").Value = 0.1323436744344018
$ws.Range("M10").Value = 134.5792873333333
$ws.Range("N10").Value = 403.737862
$ws.Range("O10").Value = 0.6702564884986638
$ws.Range("P10").Value = 0.6702564884986637
$ws.Range("Q10").Value = 416.4625130564164
$ws.Range("R10").Value = 3748.162617507748
$ws.Range("S10").Value = 0.08870420650141254
$ws.Range("T10").Value = 0.08870420650141252
